$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where the "covered" indicator (column C) becomes 1
$coveredRows = @(2, 3, 7, 8, 9, 12, 13, 25, 30, 31, 32, 33, 35, 45, 47)

foreach ($r in $coveredRows) {
    $ws.Cells.Item($r, 3).Value = 1
}

# Header formula for column C, mirroring column B's summary formula
$ws.Range("C1").Formula = "=SUM(C2:C50)/49"
$ws.Range("C1").Style = $ws.Range("B1").Style

# Row 32 (multiplyStatement) moved from 0% to 50% coverage
$ws.Range("B32").Value = 0.5
$ws.Range("A32").Style = $ws.Range("A2").Style

# Update the view: scrolled down with a new selection
$ws.Application.Goto($ws.Range("E35"), $true)
$ws.Range("E35").Select()

$wb.Save()
